$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the four test-data cells in row 5 with the refreshed values
$ws.Range("A5").Value = "TestAutomation_8Septt"
$ws.Range("B5").Value = "A224119933241"
$ws.Range("C5").Value = "Facility_h224933552q"
$ws.Range("D5").Value = "h224933552q"

# Update the active cell selection to D9
$ws.Range("D9").Select()
